$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "achievementName" -> "name"
$ws.Range("B1").Value = "name"

# Update the active selection to F12 as recorded in the saved view state
$ws.Range("F12").Select()
